$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 22.48411053931328
$ws.Range("C2").Value = 12.12474253386311
$ws.Range("D2").Value = 12.87242235501401
$ws.Range("E2").Value = 12.01070280666081
$ws.Range("G2").Value = 3.887492531174642
$ws.Range("I2").Value = 61.20418227962203
$ws.Range("J2").Value = 7.495427481341818
$ws.Range("K2").Value = 21.30710668470582
$ws.Range("L2").Value = 15.16734931949151
$ws.Range("B3").Value = 22.59371181254442
$ws.Range("C3").Value = 11.92906728949868
$ws.Range("D3").Value = 12.88698142457251
$ws.Range("E3").Value = 12.04527495031158
$ws.Range("G3").Value = 3.891799353551133
$ws.Range("I3").Value = 60.08628510378195
$ws.Range("J3").Value = 7.502022298370764
$ws.Range("K3").Value = 21.35392906340678
$ws.Range("L3").Value = 15.17499305833283
$ws.Range("B4").Value = 22.66889880102869
$ws.Range("C4").Value = 11.80937440657771
$ws.Range("D4").Value = 12.8983251591707
$ws.Range("E4").Value = 12.06805272618521
$ws.Range("G4").Value = 3.894575297762941
$ws.Range("I4").Value = 59.38666125117727
$ws.Range("J4").Value = 7.506305255564684
$ws.Range("K4").Value = 21.38882339114314
$ws.Range("L4").Value = 15.18190412831808
$ws.Range("B5").Value = 22.70151270148482
$ws.Range("C5").Value = 11.76076542501966
$ws.Range("D5").Value = 12.90355183770694
$ws.Range("E5").Value = 12.07772531922656
$ws.Range("G5").Value = 3.89573973764732
$ws.Range("I5").Value = 59.09844415897171
$ws.Range("J5").Value = 7.508109568495643
$ws.Range("K5").Value = 21.40458309579464
$ws.Range("L5").Value = 15.1852780813473
$ws.Range("B6").Value = 22.70704718674945
$ws.Range("C6").Value = 11.75270556829371
$ws.Range("D6").Value = 12.90445618458049
$ws.Range("E6").Value = 12.07935504732114
$ws.Range("G6").Value = 3.895935102629092
$ws.Range("I6").Value = 59.05040414801791
$ws.Range("J6").Value = 7.508412741104088
$ws.Range("K6").Value = 21.40729283039027
$ws.Range("L6").Value = 15.18587199973509
$ws.Range("B7").Value = 22.66933066166354
$ws.Range("C7").Value = 11.80871810314507
$ws.Range("D7").Value = 12.89839320324848
$ws.Range("E7").Value = 12.06818159237187
$ws.Range("G7").Value = 3.894590867112652
$ws.Range("I7").Value = 59.38278658061653
$ws.Range("J7").Value = 7.506329350115561
$ws.Range("K7").Value = 21.38902970415226
$ws.Range("L7").Value = 15.18194737294536
$ws.Range("B8").Value = 22.52025739052491
$ws.Range("C8").Value = 12.05721157370699
$ws.Range("D8").Value = 12.87694296946927
$ws.Range("E8").Value = 12.02230187965434
$ws.Range("G8").Value = 3.888950309691135
$ws.Range("I8").Value = 60.82160234635716
$ws.Range("J8").Value = 7.497652990121972
$ws.Range("K8").Value = 21.32197256916447
$ws.Range("L8").Value = 15.16952458614384
$ws.Range("B9").Value = 22.2909689304992
$ws.Range("C9").Value = 12.54578294919182
$ws.Range("D9").Value = 12.85397933904996
$ws.Range("E9").Value = 11.94460572883303
$ws.Range("G9").Value = 3.87892631133507
$ws.Range("I9").Value = 63.52922325953856
$ws.Range("J9").Value = 7.482483834886984
$ws.Range("K9").Value = 21.23944985248574
$ws.Range("L9").Value = 15.16275978330394
$ws.Range("B10").Value = 22.16152002836192
$ws.Range("C10").Value = 12.90246905908799
$ws.Range("D10").Value = 12.84877626258102
$ws.Range("E10").Value = 11.89496761212062
$ws.Range("G10").Value = 3.872184705621169
$ws.Range("I10").Value = 65.43777594361731
$ws.Range("J10").Value = 7.472451065109648
$ws.Range("K10").Value = 21.20894489323768
$ws.Range("L10").Value = 15.16851148427176
$ws.Range("B11").Value = 22.11121964913179
$ws.Range("C11").Value = 13.06362062828773
$ws.Range("D11").Value = 12.84894507044995
$ws.Range("E11").Value = 11.87399489102207
$ws.Range("G11").Value = 3.869251076678937
$ws.Range("I11").Value = 66.28628814356871
$ws.Range("J11").Value = 7.468125620658182
$ws.Range("K11").Value = 21.20165476032636
$ws.Range("L11").Value = 15.17345247867243
$ws.Range("B12").Value = 22.09341628628425
$ws.Range("C12").Value = 13.12443561276696
$ws.Range("D12").Value = 12.84937347482108
$ws.Range("E12").Value = 11.8662837038226
$ws.Range("G12").Value = 3.868159182485225
$ws.Range("I12").Value = 66.6045931623481
$ws.Range("J12").Value = 7.466521773031251
$ws.Range("K12").Value = 21.19984411225524
$ws.Range("L12").Value = 15.17565714173966
$ws.Range("B13").Value = 22.09719506961158
$ws.Range("C13").Value = 13.11134812260742
$ws.Range("D13").Value = 12.84926500413686
$ws.Range("E13").Value = 11.86793419196613
$ws.Range("G13").Value = 3.868393498320987
$ws.Range("I13").Value = 66.53617695286576
$ws.Range("J13").Value = 7.466865676588435
$ws.Range("K13").Value = 21.20019177954689
$ws.Range("L13").Value = 15.17516750231528
$ws.Range("B14").Value = 22.10972996453604
$ws.Range("C14").Value = 13.06862839660594
$ws.Range("D14").Value = 12.84897301227787
$ws.Range("E14").Value = 11.87335586565573
$ws.Range("G14").Value = 3.869160865757628
$ws.Range("I14").Value = 66.31253654491583
$ws.Range("J14").Value = 7.467992988754024
$ws.Range("K14").Value = 21.20148674459822
$ws.Range("L14").Value = 15.17362717570123
$ws.Range("B15").Value = 22.11757027235761
$ws.Range("C15").Value = 13.04243256901653
$ws.Range("D15").Value = 12.84884161754324
$ws.Range("E15").Value = 11.87670683050579
$ws.Range("G15").Value = 3.86963337161128
$ws.Range("I15").Value = 66.17515355765188
$ws.Range("J15").Value = 7.468687935661192
$ws.Range("K15").Value = 21.20240373454427
$ws.Range("L15").Value = 15.17272710608371
$ws.Range("B16").Value = 22.16498088466631
$ws.Range("C16").Value = 12.89191079515195
$ws.Range("D16").Value = 12.84881624650704
$ws.Range("E16").Value = 11.89637054202662
$ws.Range("G16").Value = 3.872379092662496
$ws.Range("I16").Value = 65.38191205266756
$ws.Range("J16").Value = 7.472738526003553
$ws.Range("K16").Value = 21.20955411599441
$ws.Range("L16").Value = 15.16823530193459
$ws.Range("B17").Value = 22.1962718880725
$ws.Range("C17").Value = 12.79925057072227
$ws.Range("D17").Value = 12.8494500848548
$ws.Range("E17").Value = 11.90884506530849
$ws.Range("G17").Value = 3.874097509052209
$ws.Range("I17").Value = 64.89011070273692
$ws.Range("J17").Value = 7.475284382889678
$ws.Range("K17").Value = 21.21562985062571
$ws.Range("L17").Value = 15.16607467631907
$ws.Range("B18").Value = 22.21507730224179
$ws.Range("C18").Value = 12.74585359060044
$ws.Range("D18").Value = 12.85005334514398
$ws.Range("E18").Value = 11.91617145358285
$ws.Range("G18").Value = 3.87509843891486
$ws.Range("I18").Value = 64.60539534795774
$ws.Range("J18").Value = 7.476771153923885
$ws.Range("K18").Value = 21.21974429720743
$ws.Range("L18").Value = 15.16505078420695
$ws.Range("B19").Value = 22.2215829116163
$ws.Range("C19").Value = 12.72775846756162
$ws.Range("D19").Value = 12.85029859633925
$ws.Range("E19").Value = 11.91867805613975
$ws.Range("G19").Value = 3.875439494998068
$ws.Range("I19").Value = 64.50868439516387
$ws.Range("J19").Value = 7.47727841273714
$ws.Range("K19").Value = 21.22124374150671
$ws.Range("L19").Value = 15.16474171564657
$ws.Range("B20").Value = 22.19285724772692
$ws.Range("C20").Value = 12.80912526276141
$ws.Range("D20").Value = 12.84935790840532
$ws.Range("E20").Value = 11.90750146811009
$ws.Range("G20").Value = 3.873913283645754
$ws.Range("I20").Value = 64.94265597160839
$ws.Range("J20").Value = 7.475011048804417
$ws.Range("K20").Value = 21.2149189044902
$ws.Range("L20").Value = 15.16628203345271
$ws.Range("B21").Value = 22.10601431750975
$ws.Range("C21").Value = 13.08118229684002
$ws.Range("D21").Value = 12.84904888752343
$ws.Range("E21").Value = 11.87175713064083
$ws.Range("G21").Value = 3.868934956607706
$ws.Range("I21").Value = 66.37830814098982
$ws.Range("J21").Value = 7.467660945801923
$ws.Range("K21").Value = 21.2010805808369
$ws.Range("L21").Value = 15.17407055861453
$ws.Range("B22").Value = 22.05651437361765
$ws.Range("C22").Value = 13.25774354113702
$ws.Range("D22").Value = 12.85097120503119
$ws.Range("E22").Value = 11.84974077465717
$ws.Range("G22").Value = 3.865792056436055
$ws.Range("I22").Value = 67.29898210960643
$ws.Range("J22").Value = 7.463055925306955
$ws.Range("K22").Value = 21.19757443487261
$ws.Range("L22").Value = 15.18110497407319
$ws.Range("B23").Value = 22.08226634075295
$ws.Range("C23").Value = 13.16363942856108
$ws.Range("D23").Value = 12.84975095033503
$ws.Range("E23").Value = 11.8613684426029
$ws.Range("G23").Value = 3.867459396712092
$ws.Range("I23").Value = 66.8092674767704
$ws.Range("J23").Value = 7.46549559440637
$ws.Range("K23").Value = 21.19893823682596
$ws.Range("L23").Value = 15.17717293463096
$ws.Range("B24").Value = 22.19439846632192
$ws.Range("C24").Value = 12.80466130117882
$ws.Range("D24").Value = 12.84939883737879
$ws.Range("E24").Value = 11.90810842676396
$ws.Range("G24").Value = 3.873996531482862
$ws.Range("I24").Value = 64.91890638510303
$ws.Range("J24").Value = 7.475134551103682
$ws.Range("K24").Value = 21.21523838785487
$ws.Range("L24").Value = 15.16618760735058
$ws.Range("B25").Value = 22.34618788179274
$ws.Range("C25").Value = 12.41379333168429
$ws.Range("D25").Value = 12.85814366758359
$ws.Range("E25").Value = 11.9643144815985
$ws.Range("G25").Value = 3.88152799226588
$ws.Range("I25").Value = 62.8101354467573
$ws.Range("J25").Value = 7.486391287498408
$ws.Range("K25").Value = 21.25650146105404
$ws.Range("L25").Value = 15.16270604073964
